$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D13").Value = 91.58
$ws1.Range("E13").Value = 647.08
$ws1.Range("D29").Value = "1 de 27"
$ws1.Range("E29").Value = "1 de 27"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 738.66
$ws2.Range("F29").Value = 2855.39

# --- Sheet 3: "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths: OOXML width = ColumnWidth + 0.83 (per observed conversion)
$ws3.Columns.Item(5).ColumnWidth = 23.17
$ws3.Columns.Item(6).ColumnWidth = 24.17

$ws3.Range("D3").Value = 91.58
$ws3.Range("E3").Value = 3028.5345
$ws3.Range("F3").Value = 0.02935148694062349

$ws3.Range("D4").Value = 1207.11
$ws3.Range("E4").Value = -956.4781745790989
$ws3.Range("F4").Value = 4.81626783818387

$ws3.Range("D19").Value = 6111.1
$ws3.Range("E19").Value = 17388.90093005039
$ws3.Range("F19").Value = 0.260046798218867
